# Implemented Profile section validation
# Adds a new cell (B3) on the "Profile" sheet holding the OrangeHRM support
# e-mail address, styled like the other data cells but left-aligned.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Profile")

# Seed B3 with the same look-and-feel as the existing cells (font, etc.)
# by copying the format from A1, then set its own text and alignment.
$ws.Range("A1").Copy($ws.Range("B3"))
$ws.Range("B3").Value = " ossupport@orangehrm.com "
$ws.Range("B3").HorizontalAlignment = -4131

Write-Host "Profile section validation cell added (B3)."
